$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the "currency" and "abbreviation" sheets; their single data values
# are folded into the "metadata" sheet as extra columns instead.
$currency = $wb.Worksheets.Item("currency")
[void]$currency.Delete()
$abbreviation = $wb.Worksheets.Item("abbreviation")
[void]$abbreviation.Delete()

# Rebuild "metadata" as a 2-row x 3-col table: header/currency/abbreviation
# on row 1, and the corresponding values on row 2.
$meta = $wb.Worksheets.Item("metadata")
$meta.Range("A1").Value = "header"
$meta.Range("B1").Value = "currency"
$meta.Range("C1").Value = "abbreviation"
$meta.Range("A2").Value = "Test"
$meta.Range("B2").Value = "€"
$meta.Range("C2").Value = "B"

# Match the page setup (paper size / orientation) recorded for "metadata".
$meta.PageSetup.PaperSize = 9
$meta.PageSetup.Orientation = 1

# Update the selection on "nodes": select column C, no longer the active tab.
$nodes = $wb.Worksheets.Item("nodes")
[void]$nodes.Columns.Item(3).Select()

# Update the selection on "links".
$links = $wb.Worksheets.Item("links")
[void]$links.Range("F7").Select()

# Make "metadata" the active sheet/tab and set its selection.
[void]$meta.Activate()
[void]$meta.Range("C3").Select()
